{"js": "// Office.js (Word JavaScript API) edit script.\n// Applies the same change described by the OOXML diff:\n//  1. Paragraph \" Ap\u00f3s todo esse primeiro ajuste, ... propor  melhores solu\u00e7\u00f5es. \"\n//     loses the stray double-space / grammar markup around \"propor  melhores\"\n//     becoming a clean \"...para propor melhores solu\u00e7\u00f5es. \" split over three runs.\n//  2. The first of the two trailing empty paragraphs gets the new\n//     \"3. Poderia ser aplicado o SCRUM. ...\" text.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Locate the paragraph that still contains the double-space typo\n// \"propor  melhores\" instead of searching by a fixed index (more robust\n// to any surrounding structural differences).\nlet targetParaIndex = -1;\nlet emptyParaIndex = -1;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const t = paragraphs.items[i].text;\n  if (t.indexOf(\"propor  melhores\") !== -1) {\n    targetParaIndex = i;\n  }\n}\nif (targetParaIndex !== -1) {\n  // The first empty paragraph right after the target one is where the\n  // new SCRUM paragraph needs to be written.\n  for (let i = targetParaIndex + 1; i < paragraphs.items.length; i++) {\n    if (paragraphs.items[i].text === \"\") {\n      emptyParaIndex = i;\n      break;\n    }\n  }\n}\n\nif (targetParaIndex === -1) {\n  throw new Error(\"Could not locate target paragraph containing 'propor  melhores'.\");\n}\nif (emptyParaIndex === -1) {\n  throw new Error(\"Could not locate the empty paragraph following the target paragraph.\");\n}\n\nconst targetPara = paragraphs.items[targetParaIndex];\n\n// Rebuild the paragraph's runs via OOXML so the result matches the diff\n// exactly: three runs (\"...para propor \", \"m\", \"elhores solu\u00e7\u00f5es. \") with\n// no leftover w:proofErr grammar markers (the double space that triggered\n// them is gone).\nconst wNs = 'xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"';\nconst fixedParagraphOoxml =\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  '<pkg:xmlData><w:document ' + wNs + '><w:body>' +\n  '<w:p>' +\n  '<w:r><w:t xml:space=\"preserve\"> Ap\u00f3s todo esse primeiro ajuste, chegamos na parte de execu\u00e7\u00e3o onde ser\u00e1 feito um teste semanal para ver como foi o processo de coleta seletiva na semana e assim, semanalmente ou mensalmente serem feito feedbacks para o time de planejamento para propor </w:t></w:r>' +\n  '<w:r><w:t>m</w:t></w:r>' +\n  '<w:r><w:t xml:space=\"preserve\">elhores solu\u00e7\u00f5es. </w:t></w:r>' +\n  '</w:p>' +\n  '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>';\n\ntargetPara.getRange().insertOoxml(fixedParagraphOoxml, Word.InsertLocation.replace);\nawait context.sync();\n\n// Insert the new SCRUM paragraph text into the (previously empty)\n// paragraph that followed the target paragraph.\nconst emptyPara = paragraphs.items[emptyParaIndex];\nconst scrumText =\n  \"3. Poderia ser aplicado o SCRUM. Sempre \u00e9 recomendado em qualquer projeto ter sempre um PO(Product Owner) que seria o respons\u00e1vel pela garantia do investimento e tamb\u00e9m por estar em contato direto com o cliente pra suprir suas necessidades, um SM(Scrum Master) que seria o respons\u00e1vel pela aprova\u00e7\u00e3o do projeto e da tomada de decis\u00e3o final e por fim, todo o time que ser\u00e1 feito para execu\u00e7\u00e3o do projeto.\";\nemptyPara.getRange().insertText(scrumText, Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n# Applies the same change described by the OOXML diff:\n#  1. Paragraph \" Ap\u00f3s todo esse primeiro ajuste, ... propor  melhores solu\u00e7\u00f5es. \"\n#     loses the stray double-space / grammar markup around \"propor  melhores\"\n#     becoming a clean \"...para propor melhores solu\u00e7\u00f5es. \" split over three runs.\n#  2. The first of the two trailing empty paragraphs gets the new\n#     \"3. Poderia ser aplicado o SCRUM. ...\" text.\n\n$d = $word.ActiveDocument\n\n# Locate the paragraph that still contains the double-space typo\n# \"propor  melhores\" instead of relying on a fixed index (more robust to\n# any surrounding structural differences).\n$count = $d.Paragraphs.Count\n$targetIndex = 0\nfor ($i = 1; $i -le $count; $i++) {\n    $t = $d.Paragraphs.Item($i).Range.Text\n    if ($t -like \"*propor  melhores*\") {\n        $targetIndex = $i\n    }\n}\nif ($targetIndex -eq 0) {\n    throw \"Could not locate target paragraph containing 'propor  melhores'.\"\n}\n\n# The first empty paragraph right after the target one is where the new\n# SCRUM paragraph needs to be written.\n$emptyIndex = 0\nfor ($i = $targetIndex + 1; $i -le $count; $i++) {\n    $t = $d.Paragraphs.Item($i).Range.Text\n    # A plain empty paragraph's Range.Text is just the paragraph mark.\n    if ($t -eq \"`r\" -or $t -eq \"\") {\n        $emptyIndex = $i\n        break\n    }\n}\nif ($emptyIndex -eq 0) {\n    throw \"Could not locate the empty paragraph following the target paragraph.\"\n}\n\n# Rebuild the target paragraph's runs via WordprocessingML so the result\n# matches the diff exactly: three runs (\"...para propor \", \"m\",\n# \"elhores solu\u00e7\u00f5es. \") with no leftover w:proofErr grammar markers (the\n# double space that triggered them is gone).\n$targetRange = $d.Paragraphs.Item($targetIndex).Range\n$fixedParagraphXml = '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  '<pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body>' +\n  '<w:p>' +\n  '<w:r><w:t xml:space=\"preserve\"> Ap\u00f3s todo esse primeiro ajuste, chegamos na parte de execu\u00e7\u00e3o onde ser\u00e1 feito um teste semanal para ver como foi o processo de coleta seletiva na semana e assim, semanalmente ou mensalmente serem feito feedbacks para o time de planejamento para propor </w:t></w:r>' +\n  '<w:r><w:t>m</w:t></w:r>' +\n  '<w:r><w:t xml:space=\"preserve\">elhores solu\u00e7\u00f5es. </w:t></w:r>' +\n  '</w:p>' +\n  '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'\n$targetRange.InsertXML($fixedParagraphXml)\n\n# Insert the new SCRUM paragraph text into the (previously empty)\n# paragraph that followed the target paragraph.\n$emptyRange = $d.Paragraphs.Item($emptyIndex).Range\n$emptyRange.Text = \"3. Poderia ser aplicado o SCRUM. Sempre \u00e9 recomendado em qualquer projeto ter sempre um PO(Product Owner) que seria o respons\u00e1vel pela garantia do investimento e tamb\u00e9m por estar em contato direto com o cliente pra suprir suas necessidades, um SM(Scrum Master) que seria o respons\u00e1vel pela aprova\u00e7\u00e3o do projeto e da tomada de decis\u00e3o final e por fim, todo o time que ser\u00e1 feito para execu\u00e7\u00e3o do projeto.\"\n"}
